# Fruta / hortaliza, semanal
# Insert a new weekly data row above row 12, pushing the existing
# rows 12-24 down to 13-25 (dimension grows from A1:T24 to A1:T25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 12:24 down by one row.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with this week's record.
$ws.Range("A12").Value = 7
$ws.Range("B12").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C12").Value = "Ñuble"
$ws.Range("D12").Value = 44965
$ws.Range("E12").Value = 16
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100101
$ws.Range("H12").Value = "Berries"
$ws.Range("I12").Value = 100101001
$ws.Range("J12").Value = "Arándano (blue)"
$ws.Range("K12").Value = "Sin especificar"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 50
$ws.Range("N12").Value = 3000
$ws.Range("O12").Value = 3000
$ws.Range("P12").Value = 3000
$ws.Range("Q12").Value = "$/bandeja 2 kilos"
$ws.Range("R12").Value = "Provincia de Diguillín"
$ws.Range("S12").Value = 1500
$ws.Range("T12").Value = 2
